# Generate Report for Handoff
#
# The localization run finished and is ready to hand off to the vendor:
# flip every "In Translation" status to "Ready for handoff" and bump the
# related generation timestamps. The status text is longer than before,
# so the status columns are widened to fit (mirrors an AutoFit pass on
# those columns).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Width (in Excel "characters" units) that round-trips to the widened
# stored column width used by the generated report for the status
# columns.
$statusColWidth = 16.416666666666664

# --- Overview sheet ---
# Columns E (zh-cn status) and F (de-de status).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Column G: "Latest HO Xliff Generate Date".
$wsOverview.Range("G2").Value = "2016-09-06 01:05:20"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 01:05:15"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 01:05:20"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
